$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 869 (shifts all existing rows 869-910 down to 870-911,
# which also updates the sheet's used-range dimension from D910 to D911 automatically).
$ws.Rows(869).Insert()

# Populate the newly inserted row with the new data point for 2026/02/27 (金, 13:00, rank 38).
# Force column A to be stored as literal text (matches the rest of the date column, which is
# plain text rather than a real date value) by setting the number format to "Text" before
# writing the value - otherwise Excel's COM layer auto-converts the "yyyy/mm/dd"-looking
# string into a date serial number. Reset the style back to Normal afterwards so the cell
# doesn't end up carrying a stray text-format style (the original data cells are unstyled).
$ws.Range("A869").NumberFormat = "@"
$ws.Range("A869").Value = "2026/02/27"
$ws.Range("A869").Style = "Normal"
$ws.Range("B869").Value = "金"
$ws.Range("C869").Value = 13
$ws.Range("D869").Value = 38
